# Update the "Price" column (D) of the crypto symbol list to the latest
# scraped values (GitHub Actions refresh run).
#
# The Price column stores values as text (the source sheet uses inline
# strings), so each new value is written with a leading single-quote to
# force Excel to keep it as text instead of silently re-parsing it as a
# floating point number (which would corrupt the exact decimal
# representation, e.g. "0.0006263" -> 0.0006263000000000001).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$priceUpdates = @{
    2  = "264.89"     # BNB
    3  = "22.83"      # OKB
    4  = "6.206"      # HuobiToken
    5  = "0.06159"    # Cronos
    6  = "3.563"      # GateToken
    7  = "6.703"      # KuCoinToken
    8  = "1.362"      # FTXToken
    9  = "0.8119"     # MXToken
    10 = "0.1585"     # WazirX
    11 = "0.08208"    # MandalaExchangeToken
    12 = "0.03363"    # LiechtensteinCryptoassetsExchange
    13 = "0.03149"    # BitrueCoin
    14 = "0.09246"    # BitMartToken
    15 = "3.917"      # MCDex
    16 = "0.001701"   # BitForexToken
    17 = "0.04842"    # CoinExToken
    18 = "0.0006263"  # One
    19 = "0.006174"   # TigerCash
    20 = "0.006264"   # HotbitToken
    22 = "0.0001501"  # NitroEx
    24 = "2.264"      # BTSEToken
    26 = "0.1199"     # ProBitToken
    27 = "0.0002683"  # UpBots
    40 = "0.04589"    # IDEX
    41 = "0.007004"   # KickToken
    42 = "0.1134"     # BKEXToken
    43 = "0.003132"   # CEJI
    44 = "0.01103"    # LocalTraders
    45 = "0.00006116" # CoinLion
    47 = "0.7704"     # CoinbaseStockToken
    48 = "0.1981"     # BOLO
    49 = "0.00002101" # CryptobidCoin
    50 = "0.01241"    # SpecialPowerGold
}

foreach ($row in $priceUpdates.Keys) {
    $ws.Range("D$row").Value = "'" + $priceUpdates[$row]
}
